# Adaption of Excels (remove whitespaces in column GapType!)
#
# The "Gap1_type" column (header in row 1) contains values such as
# "Arbeit", "Privat" and "Arbeit, Privat". The combined value has a
# stray space after the comma; this script strips that whitespace so
# "Arbeit, Privat" becomes "Arbeit,Privat".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

# Find the "Gap1_type" column by inspecting the header row.
$gapTypeCol = 0
for ($c = 1; $c -le $colCount; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($header -eq "Gap1_type") {
        $gapTypeCol = $c
    }
}

if ($gapTypeCol -gt 0) {
    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, $gapTypeCol)
        $val = $cell.Value()
        if ($val -ne $null -and $val.Contains(", ")) {
            $cell.Value = $val.Replace(", ", ",")
        }
    }
}
